$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.0115207373271889
$ws.Range("D2").Value = 0.0184331797235023
$ws.Range("F2").Value = 0.00691244239631336
$ws.Range("G2").Value = 0.00460829493087558
$ws.Range("H2").Value = 0.347926267281106
$ws.Range("I2").Value = 0.00921658986175115
$ws.Range("K2").Value = 0.0138248847926267
$ws.Range("L2").Value = 0.993087557603687
$ws.Range("N2").Value = 0.195852534562212
$ws.Range("Q2").Value = 0.995391705069124
$ws.Range("S2").Value = 0.914746543778802
$ws.Range("T2").Value = 0.00230414746543779
$ws.Range("U2").Value = 0.986175115207373
$ws.Range("V2").Value = 0.921658986175115
$ws.Range("W2").Value = 0.00230414746543779
$ws.Range("X2").Value = 0.00691244239631336
$ws.Range("B3").Value = 0.808755760368664
$ws.Range("D3").Value = 0.00230414746543779
$ws.Range("E3").Value = 0.00691244239631336
$ws.Range("H3").Value = 0.0115207373271889
$ws.Range("I3").Value = 0.00921658986175115
$ws.Range("J3").Value = 0.905529953917051
$ws.Range("K3").Value = 0.00460829493087558
$ws.Range("L3").Value = 0.00230414746543779
$ws.Range("M3").Value = 0.921658986175115
$ws.Range("N3").Value = 0.794930875576037
$ws.Range("P3").Value = 0.928571428571429
$ws.Range("R3").Value = 0.997695852534562
$ws.Range("S3").Value = 0.00230414746543779
$ws.Range("T3").Value = 0.995391705069124
$ws.Range("U3").Value = 0.00921658986175115
$ws.Range("V3").Value = 0.00230414746543779
$ws.Range("B4").Value = 0.0115207373271889
$ws.Range("C4").Value = 0.00460829493087558
$ws.Range("D4").Value = 0.963133640552995
$ws.Range("E4").Value = 0.00460829493087558
$ws.Range("F4").Value = 0.993087557603687
$ws.Range("G4").Value = 0.995391705069124
$ws.Range("H4").Value = 0.638248847926267
$ws.Range("J4").Value = 0.0138248847926267
$ws.Range("L4").Value = 0.00460829493087558
$ws.Range("P4").Value = 0.00460829493087558
$ws.Range("Q4").Value = 0.00230414746543779
$ws.Range("R4").Value = 0.00230414746543779
$ws.Range("S4").Value = 0.0829493087557604
$ws.Range("T4").Value = 0.00230414746543779
$ws.Range("U4").Value = 0.00460829493087558
$ws.Range("V4").Value = 0.076036866359447
$ws.Range("W4").Value = 0.995391705069124
$ws.Range("X4").Value = 0.993087557603687
$ws.Range("B5").Value = 0.179723502304147
$ws.Range("C5").Value = 0.983870967741935
$ws.Range("D5").Value = 0.0161290322580645
$ws.Range("E5").Value = 0.988479262672811
$ws.Range("H5").Value = 0.00230414746543779
$ws.Range("I5").Value = 0.981566820276498
$ws.Range("J5").Value = 0.0806451612903226
$ws.Range("K5").Value = 0.981566820276498
$ws.Range("M5").Value = 0.076036866359447
$ws.Range("N5").Value = 0.00921658986175115
$ws.Range("P5").Value = 0.0668202764976959
$ws.Range("Q5").Value = 0.00230414746543779
$ws.Range("W5").Value = 0.00230414746543779
